# Update "想去人数" (F column) values on the 展览, 演出, and 全部类型 sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15030
$ws1.Range("F3").Value = 19086
$ws1.Range("F7").Value = 222
$ws1.Range("F13").Value = 57
$ws1.Range("F15").Value = 222
$ws1.Range("F17").Value = 1473
$ws1.Range("F22").Value = 7962
$ws1.Range("F29").Value = 6063
$ws1.Range("F35").Value = 5444
$ws1.Range("F36").Value = 507

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 17

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15030
$ws4.Range("F3").Value = 19086
$ws4.Range("F7").Value = 222
$ws4.Range("F13").Value = 57
$ws4.Range("F15").Value = 222
$ws4.Range("F17").Value = 1473
$ws4.Range("F23").Value = 7962
$ws4.Range("F30").Value = 17
$ws4.Range("F32").Value = 6063
$ws4.Range("F38").Value = 5444
$ws4.Range("F39").Value = 508
